$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a "Format:" / version label row just under the title row ---
# Row 2 was an unused gap row in the original sheet (numbering jumped 1 -> 3),
# so this simply fills it in without displacing any existing rows.
$ws.Range("B2").Value = "Format:"
$ws.Range("C2").Value = "v0.1.0"

# Match the row height used throughout the rest of the sheet.
$ws.Rows.Item(2).RowHeight = 18.75

# New italic "Aptos" font for the format/version label cells.
$ws.Range("B2").Font.Italic = $true
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Name = "Aptos"
$ws.Range("B2").VerticalAlignment = -4108   # xlCenter
$ws.Range("B2").HorizontalAlignment = -4152   # xlRight

$ws.Range("C2").Font.Italic = $true
$ws.Range("C2").Font.Size = 11
$ws.Range("C2").Font.Name = "Aptos"
$ws.Range("C2").VerticalAlignment = -4108   # xlCenter
